$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> MuSCs (target cluster), updated TPM-derived metrics ---
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 82.98768099999999
$ws.Range("H2").Value = 248.963043
$ws.Range("I2").Value = 0.4489504115427952
$ws.Range("J2").Value = 0.4489504115427952
$ws.Range("M2").Value = 0.000484
$ws.Range("N2").Value = 0.001452
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.040166037604
$ws.Range("R2").Value = 0.361494338436
$ws.Range("S2").Value = 0.4489504115427952
$ws.Range("T2").Value = 0.4489504115427952

# --- Row 3: ECs -> FAPs (sending cluster), target cluster stays MuSCs ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 63.14058933333333
$ws.Range("H3").Value = 189.421768
$ws.Range("I3").Value = 0.3415807409566563
$ws.Range("J3").Value = 0.3415807409566563
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.03056004523733333
$ws.Range("R3").Value = 0.275040407136
$ws.Range("S3").Value = 0.3415807409566563
$ws.Range("T3").Value = 0.3415807409566563

# --- Row 4: FAPs -> MuSCs (sending cluster), ECs -> MuSCs (target cluster) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 38.719942
$ws.Range("H4").Value = 116.159826
$ws.Range("I4").Value = 0.2094688475005485
$ws.Range("J4").Value = 0.2094688475005485
$ws.Range("M4").Value = 0.000484
$ws.Range("N4").Value = 0.001452
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.018740451928
$ws.Range("R4").Value = 0.168664067352
$ws.Range("S4").Value = 0.2094688475005485
$ws.Range("T4").Value = 0.2094688475005485

# --- Remove the now-obsolete rows 5-7 (previous MuSCs-target duplicates) ---
$ws.Rows("5:7").Delete()
